$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on columns that would otherwise be auto-converted
# by Excel (numeric-looking strings in "Antal", and date-looking strings
# in Startdatum/Slutdatum), to preserve them as plain text like the source data.
$ws.Range("I2:I6").NumberFormat = "@"
$ws.Range("Y2:AB6").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 109910753
$ws.Range("B2").Value = 78596
$ws.Range("C2").Value = "Ovaliderad"
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 6462
$ws.Range("F2").Value = "Stuplav"
$ws.Range("G2").Value = "Nephroma bellum"
$ws.Range("H2").Value = "(Spreng.) Tuck."
$ws.Range("P2").Value = "Risberg, Vb"
$ws.Range("Q2").Value = 728699.2894223372
$ws.Range("R2").Value = 7192686.04447869
$ws.Range("S2").Value = 10
$ws.Range("T2").Value = "Västerbotten"
$ws.Range("U2").Value = "Norsjö"
$ws.Range("V2").Value = "Västerbotten"
$ws.Range("W2").Value = "Norsjö"
$ws.Range("Y2").Value = "2022-11-01"
$ws.Range("Z2").Value = "00:00"
$ws.Range("AA2").Value = "2022-11-01"
$ws.Range("AB2").Value = "00:00"
$ws.Range("AC2").Value = "Påträffad under Sveaskogs naturvärdesbedömning"
$ws.Range("AD2").Value = $false
$ws.Range("AE2").Value = $false
$ws.Range("AG2").Value = $false
$ws.Range("AW2").Value = "Mimmi Persson"
$ws.Range("AX2").Value = "Mimmi Persson"
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("K2").ClearContents()
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").ClearContents()

# Row 3
$ws.Range("A3").Value = 6803431
$ws.Range("B3").Value = 96354
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "LC"
$ws.Range("E3").Value = 221952
$ws.Range("F3").Value = "Spindelblomster"
$ws.Range("G3").Value = "Neottia cordata"
$ws.Range("H3").Value = "(L.) Rich."
$ws.Range("I3").Value = "15"
$ws.Range("J3").Value = "plantor/tuvor"
$ws.Range("K3").Value = "blomning"
$ws.Range("P3").Value = "storbäcken, Vb"
$ws.Range("Q3").Value = 729732.4712969258
$ws.Range("R3").Value = 7192285.441698131
$ws.Range("S3").Value = 50
$ws.Range("T3").Value = "Västerbotten"
$ws.Range("U3").Value = "Norsjö"
$ws.Range("V3").Value = "Västerbotten"
$ws.Range("W3").Value = "Norsjö"
$ws.Range("Y3").Value = "2013-06-18"
$ws.Range("Z3").Value = "10:00"
$ws.Range("AA3").Value = "2013-06-18"
$ws.Range("AB3").Value = "10:00"
$ws.Range("AC3").Value = "förekommer i fuktstråk i äldre granskog längs storbäcken"
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false
$ws.Range("AH3").Value = "Högörtgranskog"
$ws.Range("AI3").Value = "fuktstråk i äldre granskog omkring storbäcken"
$ws.Range("AW3").Value = "Mikael Marberg"
$ws.Range("AX3").Value = "Mikael Marberg"
$ws.Range("AJ3").ClearContents()
$ws.Range("AK3").ClearContents()
$ws.Range("AL3").ClearContents()
$ws.Range("AO3").ClearContents()

# Row 4
$ws.Range("A4").Value = 6803307
$ws.Range("B4").Value = 78569
$ws.Range("C4").Value = "Ovaliderad"
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 6458
$ws.Range("F4").Value = "Lunglav"
$ws.Range("G4").Value = "Lobaria pulmonaria"
$ws.Range("H4").Value = "(L.) Hoffm."
$ws.Range("J4").Value = "bålar"
$ws.Range("P4").Value = "storbäcken, Vb"
$ws.Range("Q4").Value = 729813.2748242758
$ws.Range("R4").Value = 7192204.280184632
$ws.Range("S4").Value = 50
$ws.Range("T4").Value = "Västerbotten"
$ws.Range("U4").Value = "Norsjö"
$ws.Range("V4").Value = "Västerbotten"
$ws.Range("W4").Value = "Norsjö"
$ws.Range("Y4").Value = "2013-06-17"
$ws.Range("Z4").Value = "15:13"
$ws.Range("AA4").Value = "2013-06-17"
$ws.Range("AB4").Value = "15:14"
$ws.Range("AC4").Value = "Förekommer på äldre sälgar söder om myren och österut längs storbäcken."
$ws.Range("AD4").Value = $false
$ws.Range("AE4").Value = $false
$ws.Range("AG4").Value = $false
$ws.Range("AJ4").Value = "vanlig sälg"
$ws.Range("AK4").Value = "Salix caprea subsp. caprea"
$ws.Range("AL4").Value = "äldre sälgar i avsatt granskog längs bäck och våtmark"
$ws.Range("AO4").Value = "Salix caprea subsp. caprea # äldre sälgar i avsatt granskog längs bäck och våtmark"
$ws.Range("AW4").Value = "Mikael Marberg"
$ws.Range("AX4").Value = "Mikael Marberg"
$ws.Range("I4").ClearContents()
$ws.Range("AI4").ClearContents()

# Row 5
$ws.Range("A5").Value = 6803316
$ws.Range("B5").Value = 89410
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("I5").Value = "15"
$ws.Range("J5").Value = "fruktkroppar"
$ws.Range("P5").Value = "storbäcken, Vb"
$ws.Range("Q5").Value = 729732.4712969258
$ws.Range("R5").Value = 7192285.441698131
$ws.Range("S5").Value = 50
$ws.Range("T5").Value = "Västerbotten"
$ws.Range("U5").Value = "Norsjö"
$ws.Range("V5").Value = "Västerbotten"
$ws.Range("W5").Value = "Norsjö"
$ws.Range("Y5").Value = "2013-06-18"
$ws.Range("Z5").Value = "11:12"
$ws.Range("AA5").Value = "2013-06-18"
$ws.Range("AB5").Value = "11:15"
$ws.Range("AC5").Value = "funnen på två granlågor"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AI5").Value = "äldre granskog med fuktstråk som avsatts längs bäck"
$ws.Range("AW5").Value = "Mikael Marberg"
$ws.Range("AX5").Value = "Mikael Marberg"
$ws.Range("K5").ClearContents()
$ws.Range("AH5").ClearContents()

# Row 6
$ws.Range("A6").Value = 6803279
$ws.Range("B6").Value = 98493
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 1365
$ws.Range("F6").Value = "Lappranunkel"
$ws.Range("G6").Value = "Coptidium lapponicum"
$ws.Range("H6").Value = "(L.) Tzvelev"
$ws.Range("I6").Value = "300"
$ws.Range("J6").Value = "stjälkar/strån/skott"
$ws.Range("K6").Value = "blomning"
$ws.Range("P6").Value = "storbäcken, Vb"
$ws.Range("Q6").Value = 729813.2748242758
$ws.Range("R6").Value = 7192204.280184632
$ws.Range("S6").Value = 50
$ws.Range("T6").Value = "Västerbotten"
$ws.Range("U6").Value = "Norsjö"
$ws.Range("V6").Value = "Västerbotten"
$ws.Range("W6").Value = "Norsjö"
$ws.Range("Y6").Value = "2013-06-18"
$ws.Range("Z6").Value = "11:00"
$ws.Range("AA6").Value = "2013-06-18"
$ws.Range("AB6").Value = "11:30"
$ws.Range("AC6").Value = "över 300 blommande"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AH6").Value = "Gransumpskog"
$ws.Range("AI6").Value = "genomsilande markvatten i äldre bäcknära granskog"
$ws.Range("AW6").Value = "Mikael Marberg"
$ws.Range("AX6").Value = "Mikael Marberg"
